$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.269.81"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").Value = "1.595.65"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'213.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.78%  "
$ws.Range("E6").Value = "  -0.81%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.246"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.34%  "
$ws.Range("E9").Value = "  -0.64%  "
$ws.Range("D10").Value = "'19.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.52%  "
$ws.Range("E11").Value = "  +0.49%  "
$ws.Range("D12").Value = "1.819.95"
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("D13").Value = "1.619.50"
$ws.Range("E13").Value = "  +1.68%  "
$ws.Range("E14").Value = "  -1.37%  "
$ws.Range("E15").Value = "  -2.69%  "
$ws.Range("D16").Value = "'63.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.38%  "
$ws.Range("D17").Value = "26.254.91"
$ws.Range("E17").Value = "  -0.44%  "
$ws.Range("E18").Value = "  -0.73%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'7.37"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.11%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'214.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.91%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").Value = "'9.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.42%  "
$ws.Range("E24").Value = "  -3.47%  "
$ws.Range("D25").Value = "'145.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.30%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "'6.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.61%  "
$ws.Range("E28").Value = "  -0.77%  "
$ws.Range("D29").Value = "'15.13"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.67%  "
$ws.Range("E30").Value = "  -2.74%  "
$ws.Range("E31").Value = "  +0.60%  "
$ws.Range("E32").Value = "  -0.67%  "
$ws.Range("D33").Value = "1.424.69"
$ws.Range("E33").Value = "  +6.43%  "
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("D35").Value = "'2.42"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.69%  "
$ws.Range("E36").Value = "  -1.02%  "
$ws.Range("D37").Value = "'0.581"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.42%  "
$ws.Range("E38").Value = "  -0.80%  "
$ws.Range("E39").Value = "  +0.64%  "
$ws.Range("E40").Value = "  +0.75%  "
$ws.Range("D42").Value = "'0.975"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.35%  "
$ws.Range("E43").Value = "  +1.16%  "
$ws.Range("D44").Value = "'0.766"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("D45").Value = "1.731.24"
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("D46").Value = "'61.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.48%  "
$ws.Range("D47").Value = "'87.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.34%  "
$ws.Range("D48").Value = "'1.50"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("E50").Value = "  -2.72%  "
$ws.Range("D51").Value = "'1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.04%  "
